$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column AI (35), shifting "nom"/"url_produit" right.
$ws.Columns.Item(35).Insert()

# New header in AI1: next timestamp in the price-history sequence.
$ws.Cells.Item(1, 35).Value = "2026-01-29 06:32:37"

# For each data row, the new AI column carries forward the last price value (AH).
$lastRow = 205
for ($r = 2; $r -le $lastRow; $r++) {
    $ahVal = $ws.Cells.Item($r, 34).Value2
    if ($null -ne $ahVal -and $ahVal -ne "") {
        $ws.Cells.Item($r, 35).Value = $ahVal
    } else {
        $ws.Cells.Item($r, 35).Value = ""
    }
}
